$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint Backlog")

$ws.Range("K4").Value = 3
$ws.Range("K5").Value = 3
$ws.Range("K7").Value = 5
$ws.Range("K8").Value = 4

$ws.Activate()
$ws.Range("K4").Select()
